$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 15/10/2013, duration 02:05 (0.0868055555555556 of a day)
# Copy formatting down from the row above so the new cells reuse the
# existing date/time number-format styles instead of creating new ones.
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)

$ws.Range("A10").Value = 41562
$ws.Range("B10").Value = 0.086805555555555566

# Move/collapse the selection to C10, matching the saved view state.
$ws.Range("C10").Select()
